# Auto-committed on 2022/08/26 週五 17:20:49.47
#
# Change the "形態" (data type) of the 段/巷/弄/號/號之/樓/樓之 address fields
# (both 戶籍/公司 "Reg*" rows and 通訊 "Curr*" rows) from VARCHAR2 to
# NVARCHAR2 so they can store Chinese text, and annotate each changed row
# with a note in a new column H explaining why.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$note = "2022-08-24 可以打中文 所以要改型態為NVARCHAR2"

# Rows 27-33 (RegSection/RegAlley/RegLane/RegNum/RegNumDash/RegFloor/RegFloorDash)
# and rows 39-45 (CurrSection/CurrAlley/CurrLane/CurrNum/CurrNumDash/CurrFloor/CurrFloorDash)
# all get their "形態" column switched from VARCHAR2 to NVARCHAR2 (with the
# matching highlighted style) and a new remark in column H.
$targetRows = @(27, 28, 29, 30, 31, 32, 33, 39, 40, 41, 42, 43, 44, 45)

foreach ($r in $targetRows) {
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = "NVARCHAR2"
    $dCell.Interior.Color = 65535

    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $note
}

# Scroll/zoom/selection bookkeeping to match the saved view state.
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F37").Select()
